# Append a new log row (row 26) to the Nalco run log, recording a SKIPPED
# run performed on 2025-08-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "2025-08-18 04:12:06 UTC"
$ws.Range("B26").Value = "2025-08-18 09:42:06 IST"
$ws.Range("C26").Value = "SKIPPED"
$ws.Range("D26").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E26").Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = ""

# Match the existing data-row formatting (centered, no border - style index 3)
# used by every other row in the sheet.
$ws.Range("A26:H26").HorizontalAlignment = -4108
$ws.Range("A26:H26").VerticalAlignment = -4108

Write-Output "Appended row 26 to $($ws.Name)"
